# Add newly tracked coding-question rows (49-54) plus two trailing index
# rows (67-68) to the "CodingQuestionsHint" tracker sheet, matching the
# upstream "Add files via upload" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Copy-Format {
    # Copies formatting only (style) from one cell to another, the same
    # way a user would with the Format Painter / Paste Special > Formats.
    param(
        [string]$FromAddr,
        [string]$ToAddr
    )
    $ws.Range($FromAddr).Copy() | Out-Null
    $ws.Range($ToAddr).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
}

function Set-LeetcodeCell {
    # Writes "<prefix>Leetcode" into a cell with the trailing "Leetcode"
    # word rendered bold (matches the source-of-truth rich text runs).
    param(
        [string]$CellAddr,
        [string]$Prefix
    )
    $bold = "Leetcode"
    $ws.Range($CellAddr).Value = $Prefix + $bold
    $ws.Range($CellAddr).Characters($Prefix.Length + 1, $bold.Length).Font.Bold = $true
}

# ----------------------------------------------------------------------
# Row 49
# ----------------------------------------------------------------------
$ws.Range("A49").Value = 47

Copy-Format "B48" "B49"
$ws.Range("B49").Value = 45661

Copy-Format "C4" "G49"
$ws.Range("G49").Value = "solved but not submitted"

Copy-Format "C4" "F49"
$ws.Range("F49").Value = "O(n^2)."

Copy-Format "C4" "C49"
Set-LeetcodeCell "C49" ".3sum. Given an integer array nums, return all the triplets [nums[i], nums[j], nums[k]] such that i != j, i != k, and j != k, and nums[i] + nums[j] + nums[k] == 0.. "

Copy-Format "D48" "D49"
$ws.Range("D49").Value = " [-1,0,1,2,-1,-4]"

Copy-Format "E48" "E49"
$ws.Range("E49").Value = "[[-1,-1,2],[-1,0,1]]"

# ----------------------------------------------------------------------
# Row 50
# ----------------------------------------------------------------------
$ws.Range("A50").Value = 48

Copy-Format "B49" "B50"
$ws.Range("B50").Value = 45662

Copy-Format "C4" "C50"
Set-LeetcodeCell "C50" "find index of first occurrence in string. "

Copy-Format "C4" "D50"
$ws.Range("D50").Value = 'haystack = "sadbutsad", needle = "sad"'

$ws.Range("E50").Value = 0

Copy-Format "C4" "G50"
$ws.Range("G50").Value = "took help"

# ----------------------------------------------------------------------
# Row 51
# ----------------------------------------------------------------------
$ws.Range("A51").Value = 49

Copy-Format "B50" "B51"
$ws.Range("B51").Value = 45662

Copy-Format "C4" "C51"
Set-LeetcodeCell "C51" "Binary tree inorder traversal. "

Copy-Format "C4" "G51"
$ws.Range("G51").Value = "solved but not submitted"

# ----------------------------------------------------------------------
# Row 52
# ----------------------------------------------------------------------
$ws.Range("A52").Value = 50

Copy-Format "B51" "B52"
$ws.Range("B52").Value = 45662

Copy-Format "C4" "C52"
Set-LeetcodeCell "C52" "Binary tree preorder traversal. "

Copy-Format "C4" "G52"
$ws.Range("G52").Value = "solved but not submitted"

# ----------------------------------------------------------------------
# Row 53
# ----------------------------------------------------------------------
$ws.Range("A53").Value = 51

Copy-Format "B52" "B53"
$ws.Range("B53").Value = 45662

Copy-Format "C4" "C53"
Set-LeetcodeCell "C53" "Binary tree postorder traversal. "

Copy-Format "C4" "G53"
$ws.Range("G53").Value = "solved but not submitted"

# ----------------------------------------------------------------------
# Row 54
# ----------------------------------------------------------------------
$ws.Range("A54").Value = 52

Copy-Format "B53" "B54"
$ws.Range("B54").Value = 45662

Copy-Format "C4" "C54"
Set-LeetcodeCell "C54" "Seach in sorted rotated array. "

Copy-Format "C4" "F54"
$ws.Range("F54").Value = "O(logn)"

Copy-Format "C4" "G54"
$ws.Range("G54").Value = "solved but not submitted"

# ----------------------------------------------------------------------
# Two additional trailing index rows (67, 68), continuing the numbering
# already running down column A.
# ----------------------------------------------------------------------
$ws.Range("A67").Value = 65
$ws.Range("A68").Value = 66

# ----------------------------------------------------------------------
# Update the view to where the author left off editing.
# ----------------------------------------------------------------------
$ws.Range("G54").Select() | Out-Null
